# Adds a "total" column (BS) with row sums, plus a grand-total cell at the
# bottom of that column, to the concept-relationship matrix on Sheet1.
#
# Data lives in A1:BR70:
#   - Row 1 / Column A hold the same list of concept labels.
#   - B2:BR70 is a 0/1 adjacency-style matrix.
# New additions:
#   - BS2:BS70 => =SUM(B<r>:BR<r>)  (row total for each concept)
#   - BS71     => =SUM(BS2:BS70)    (grand total)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 70
$totalCol = "BS"

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Range("$totalCol$r").Formula = "=SUM(B$r`:BR$r)"
}

$grandTotalRow = $lastDataRow + 1
$ws.Range("$totalCol$grandTotalRow").Formula = "=SUM($totalCol$firstDataRow`:$totalCol$lastDataRow)"

# Reflect the final selection/scroll state used when the workbook was last
# saved: the bottom total row selected, scrolled near the bottom of the data.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 45
$win.ScrollColumn = 1
$ws.Range("B70:BR70").Select() | Out-Null
